$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New response row (row 6) - "Zain" entry being added to the form responses sheet.

# Apply number formats to the new row first (matching the existing rows' formats)
# so that subsequently-written values land in correctly-formatted cells.
$ws.Range("A6").NumberFormat = "m/d/yyyy h:mm:ss"
$ws.Range("C6:AZ6").NumberFormat = "General"
$ws.Range("BA6").NumberFormat = "General"

# Timestamp
$ws.Range("A6").Value = 44550.538091504626

# Name
$ws.Range("B6").Value = "Zain"

# Ranking columns C..AZ (50 columns of rank values, matching form question order)
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
          "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT", `
          "AU","AV","AW","AX","AY","AZ")
$vals = @(3,1,2,5,4,7,10,8,6,9,2,3,5,4,1,6,9,10,7,8,7,4,2,1,6,5,9,8,3,10,5,3,1,2,4,8,9,10,6,7,5,6,1,2,3,7,4,10,8,9)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $vals[$i]
}

# Score column
$ws.Range("BA6").Value = 0.0
